$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $range = $ws.Range($ref)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# total changes: 80
Set-TextCell "D2" "26.716.64"
Set-TextCell "D3" "1.600.77"
Set-TextCell "E3" "  +0.21%  "
Set-TextCell "E4" "  +0.29%  "
Set-TextCell "D5" "211.51"
Set-TextCell "E5" "  -0.13%  "
Set-TextCell "E6" "  -0.69%  "
Set-TextCell "E7" "  +0.28%  "
Set-TextCell "E8" "  +0.14%  "
Set-TextCell "E9" "  +0.72%  "
Set-TextCell "E10" "  +0.27%  "
Set-TextCell "D11" "0.0843"
Set-TextCell "E11" "  +0.64%  "
Set-TextCell "D12" "1.825.17"
Set-TextCell "D13" "1.590.74"
Set-TextCell "E13" "  -0.46%  "
Set-TextCell "E14" "  +0.50%  "
Set-TextCell "E15" "  +0.14%  "
Set-TextCell "D16" "65.37"
Set-TextCell "E16" "  +1.38%  "
Set-TextCell "D17" "26.690.80"
Set-TextCell "E17" "  +0.21%  "
Set-TextCell "E18" "  +3.00%  "
Set-TextCell "E19" "  +3.65%  "
Set-TextCell "E20" "  +0.29%  "
Set-TextCell "D21" "209.12"
Set-TextCell "E21" "  -0.02%  "
Set-TextCell "D22" "4.30"
Set-TextCell "E22" "  +0.54%  "
Set-TextCell "E23" "  +0.86%  "
Set-TextCell "E24" "  +0.65%  "
Set-TextCell "D25" "142.54"
Set-TextCell "E25" "  -1.91%  "
Set-TextCell "D26" "1.00"
Set-TextCell "E26" "  +0.15%  "
Set-TextCell "D27" "7.11"
Set-TextCell "E27" "  -0.75%  "
Set-TextCell "E28" "  +0.07%  "
Set-TextCell "E29" "  +0.52%  "
Set-TextCell "E30" "  +2.95%  "
Set-TextCell "D31" "1.15"
Set-TextCell "E31" "  -0.32%  "
Set-TextCell "D32" "3.25"
Set-TextCell "E32" "  +0.54%  "
Set-TextCell "E33" "  +1.56%  "
Set-TextCell "D34" "1.293.62"
Set-TextCell "E34" "  +1.35%  "
Set-TextCell "E35" "  -5.03%  "
Set-TextCell "E36" "  +0.96%  "
Set-TextCell "E37" "  +0.58%  "
Set-TextCell "E38" "  -0.26%  "
Set-TextCell "E39" "  +20.24%  "
Set-TextCell "E41" "  -0.81%  "
Set-TextCell "D42" "2.21"
Set-TextCell "E42" "  +0.34%  "
Set-TextCell "D43" "0.783"
Set-TextCell "E43" "  -0.39%  "
Set-TextCell "D44" "63.15"
Set-TextCell "E44" "  -2.28%  "
Set-TextCell "D45" "1.736.88"
Set-TextCell "E45" "  +0.12%  "
Set-TextCell "D46" "91.37"
Set-TextCell "E46" "  +1.48%  "
Set-TextCell "E47" "  -1.92%  "
Set-TextCell "B48" "BabyDogeCoin"
Set-TextCell "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell "D48" "0.0₆0105"
Set-TextCell "E48" "  -0.74%  "
Set-TextCell "B49" "Algorand"
Set-TextCell "C49" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextCell "D49" "0.101"
Set-TextCell "E49" "  -1.36%  "
Set-TextCell "B50" "Cronos"
Set-TextCell "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D50" "0.0511"
Set-TextCell "E50" "  +0.58%  "
Set-TextCell "B51" "USDD"
Set-TextCell "C51" "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextCell "D51" "1.00"
Set-TextCell "E51" "  +0.21%  "
